# Logged Week 16 and performed season sim from Week 17
# Update the "R" (Road) row target-depth splits on both the OFF and DEF
# sheets with the latest simulated/logged totals.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 is the "R" row ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 188
$wsOff.Range("C3").Value = 138
$wsOff.Range("D3").Value = 45
$wsOff.Range("E3").Value = 19

# --- DEF sheet: row 3 is the "R" row ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 178
$wsDef.Range("C3").Value = 110
$wsDef.Range("D3").Value = 47
$wsDef.Range("E3").Value = 24
$wsDef.Range("F3").Value = 4
